# Update the footer/date placeholder's cached "datetimeFigureOut" field text
# from 2020/7/19 to 2020/8/9 across the slide master and every slide layout.

$p = $ppt.ActivePresentation
$ppPlaceholderDate = 16
$oldDate = "2020/7/19"
$newDate = "2020/8/9"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master.
$master = $p.SlideMaster
Update-DateShape $master.Shapes

# Every slide layout (CustomLayout) hanging off the master.
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    Update-DateShape $layout.Shapes
}
